$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price column cells to stay as text so strings like
# "509.50" / "7.40" / "0.850" are not auto-converted to numbers
# (which would silently drop the significant trailing zero).
$priceCells = @("D2","D3","D5","D6","D9","D10","D14","D15","D16","D18","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D32","D33","D35","D36","D40","D41","D42","D43","D44","D45","D46","D48","D49","D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "60.651.73"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("D3").Value = "2.591.46"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "509.50"
$ws.Range("E5").Value = "  -1.81%  "
$ws.Range("D6").Value = "155.75"
$ws.Range("E6").Value = "  -3.35%  "
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("E8").Value = "  -4.85%  "
$ws.Range("D9").Value = "2.599.65"
$ws.Range("E9").Value = "  -3.39%  "
$ws.Range("D10").Value = "6.56"
$ws.Range("E10").Value = "  +6.27%  "
$ws.Range("E11").Value = "  -3.12%  "
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("D14").Value = "3.042.67"
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("D15").Value = "60.592.62"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").Value = "21.71"
$ws.Range("E16").Value = "  -4.08%  "
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("D18").Value = "2.595.20"
$ws.Range("E18").Value = "  -2.55%  "
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("D20").Value = "347.22"
$ws.Range("E20").Value = "  -4.64%  "
$ws.Range("D21").Value = "10.53"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D22").Value = "6.13"
$ws.Range("E22").Value = "  -1.92%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "60.26"
$ws.Range("E24").Value = "  -1.35%  "
$ws.Range("D25").Value = "0.421"
$ws.Range("E25").Value = "  -2.00%  "
$ws.Range("D26").Value = "0.168"
$ws.Range("E26").Value = "  -1.86%  "
$ws.Range("D27").Value = "2.702.12"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("D29").Value = "0.0₃0850"
$ws.Range("E29").Value = "  -3.47%  "
$ws.Range("D30").Value = "7.40"
$ws.Range("E30").Value = "  -3.37%  "
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("D32").Value = "19.45"
$ws.Range("E32").Value = "  -2.45%  "
$ws.Range("D33").Value = "152.88"
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("E34").Value = "  -2.05%  "
$ws.Range("D35").Value = "5.72"
$ws.Range("E35").Value = "  +1.02%  "
$ws.Range("D36").Value = "4.02"
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("E37").Value = "  -3.27%  "
$ws.Range("E38").Value = "  +3.47%  "
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("D40").Value = "0.850"
$ws.Range("E40").Value = "  -4.14%  "
$ws.Range("D41").Value = "36.22"
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("D42").Value = "3.76"
$ws.Range("E42").Value = "  -1.23%  "
$ws.Range("D43").Value = "297.95"
$ws.Range("E43").Value = "  -4.30%  "
$ws.Range("D44").Value = "0.625"
$ws.Range("E44").Value = "  -3.48%  "
$ws.Range("D45").Value = "0.0996"
$ws.Range("E45").Value = "  -2.67%  "
$ws.Range("D46").Value = "0.0561"
$ws.Range("E46").Value = "  -3.63%  "
$ws.Range("E47").Value = "  +1.39%  "
$ws.Range("D48").Value = "19.82"
$ws.Range("E48").Value = "  -1.81%  "
$ws.Range("D49").Value = "4.86"
$ws.Range("E49").Value = "  -4.21%  "
$ws.Range("D50").Value = "0.0234"
$ws.Range("E50").Value = "  -2.86%  "
$ws.Range("E51").Value = "  +0.21%  "

# Restore the default (general) style on the price cells so we
# do not leave a stray explicit number-format style behind.
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
